# "Generate Report for Handoff"
#
# The handback-status for the e2e\9dced597-287e-4541-a4e2-02146d0d2a6d.md
# file moves from "Handed back: in sync with en-US" to "Ready for handoff",
# with refreshed timestamps and a new "version not latest" error message,
# on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2392831e807ccb6fbb269c95a44cec6be411fb52/e2e/9dced597-287e-4541-a4e2-02146d0d2a6d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f021f8b83fff13ccdd9d275ac98ad90208141a47/e2e/9dced597-287e-4541-a4e2-02146d0d2a6d.md."

# ---- Overview sheet: row 3 is the 9dced597...md file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 09:02:26"

# ---- zh-cn sheet: row 3 is the 9dced597...md file ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-01 09:02:22"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1:P3").ColumnWidth = 39.14

# ---- de-de sheet: row 3 is the 9dced597...md file ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-01 09:02:26"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1:P3").ColumnWidth = 39.14
